$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume values (and the WrappedBTC/WrappedEther row swap).
# A leading apostrophe forces Excel to keep the assigned text as a literal string
# (rather than re-interpreting numeric-looking text such as "706.33" as a float),
# and resetting the Style back to "Normal" drops the temporary quote-prefix formatting
# so the cell ends up with the same (default) style as before.
$ws.Range('D2').Value = '''71.222.32'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +0.49%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''3.813.99'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -0.80%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.01%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''706.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +1.42%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''171.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +0.06%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''3.812.89'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -0.79%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.01%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.22%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.161'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -0.12%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''7.68'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +5.63%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.462'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +0.83%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = '''  -1.40%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''36.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -0.26%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''4.459.90'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -0.70%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = '''WrappedEther'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value = '''3.853.46'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +0.34%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = '''WrappedBTC'
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value = '''https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value = '''71.219.07'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  +0.49%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''17.55'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +1.30%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''7.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +0.24%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = '''  -0.23%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''514.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +4.24%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''10.70'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -0.18%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = '''  +0.86%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''84.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -0.44%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = '''  -1.67%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''3.967.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -0.81%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '''  -1.05%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''10.43'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -1.13%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = '''  +0.17%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = '''  -3.28%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = '''  -2.52%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''7.42'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -0.89%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = '''  -0.36%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''29.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -1.03%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  -3.90%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''9.18'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +0.57%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''3.777.64'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -0.62%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = '''  -0.08%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = '''  -1.86%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''2.38'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -0.30%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''6.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +0.25%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''1.03'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  -1.37%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = '''  -1.58%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = '''  -0.02%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''171.10'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +4.68%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = '''  +0.14%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''0.000311'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  +0.27%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''49.63'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +1.78%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''423.03'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  +4.33%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  +0.39%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.294'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -1.49%  '
$ws.Range('E51').Style = 'Normal'
